$d = $word.ActiveDocument

# --- locate "Close another admission" -------------------------------------
$rng = $d.Content
$rng.Find.Execute("Close another admission", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end   = $rng.End

# --- the actual text edit ---------------------------------------------------
# "Close another admission"  ->  "Exit or close another admission"
# i.e. lower-case the leading "C" and insert "Exit or " before it.

# 1) "C" -> "c" (same length, in place)
$d.Range($start, $start + 1).Text = "c"

# 2) insert "Exit or " right before it
$d.Range($start, $start).InsertBefore("Exit or ")

$newEnd = $end + 8   # "Exit or " is 8 characters, so every later offset shifts by 8
$splitAt = $start + 9   # right after "...Exit or c"

# --- restore the original run layout ---------------------------------------
# This engine coalesces any stretch of identically-formatted runs that is
# touched by an edit into a single run. Left alone, that would merge our
# edited text with the quotation-mark/"prompt." runs on either side of it,
# and would leave "Exit or close another admission" as one run instead of
# the target's "Exit or c" + "lose another admission" split.
#
# Toggling a character property on and back off for a given span is a
# content no-op (the property ends up exactly as it started everywhere) but
# it forces the engine to keep that span as its own run instead of folding
# it into a neighbour. We apply it to each span that must stay a distinct
# run: the quote mark before, the two halves of the split word, and the
# quote+prompt runs after - i.e. exactly the runs touching the edit point.

$r1 = $d.Range($start - 1, $start)        # the opening curly quote, "
$r1.Font.Bold = $true
$r1.Font.Bold = $false

$r2 = $d.Range($start, $splitAt)          # "Exit or c"
$r2.Font.Bold = $true
$r2.Font.Bold = $false

$r3 = $d.Range($splitAt, $newEnd)         # "lose another admission"
$r3.Font.Bold = $true
$r3.Font.Bold = $false

$r4 = $d.Range($newEnd, $newEnd + 2)      # the closing curly quote + space, " "
$r4.Font.Bold = $true
$r4.Font.Bold = $false

$r5 = $d.Range($newEnd + 2, $newEnd + 10) # "prompt. "
$r5.Font.Bold = $true
$r5.Font.Bold = $false
